$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.239.71'
$ws.Range("E2").Value = '  +0.57%  '

$ws.Range("D3").Value = '2.586.61'
$ws.Range("E3").Value = '  +2.14%  '

$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '316.04'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.43%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '97.00'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.26%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.578'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.17%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.05%  '

$ws.Range("E9").Value = '  +1.30%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.70'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.29%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0815'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.46%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.51'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.09%  '

$ws.Range("D13").Value = '2.978.04'
$ws.Range("E13").Value = '  +1.96%  '

$ws.Range("E14").Value = '  -3.23%  '

$ws.Range("D15").Value = '2.500.98'
$ws.Range("E15").Value = '  -3.44%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.20'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.06%  '

$ws.Range("E17").Value = '  -0.58%  '

$ws.Range("D18").Value = '43.255.25'
$ws.Range("E18").Value = '  +0.57%  '

$ws.Range("E19").Value = '  +2.86%  '

$ws.Range("E20").Value = '  -3.51%  '

$ws.Range("D21").Value = '0.0₃0964'
$ws.Range("E21").Value = '  -0.28%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '69.53'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.02%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '254.51'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.05%  '

$ws.Range("E24").Value = '  +0.74%  '

$ws.Range("E25").Value = '  +2.84%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '27.38'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.52%  '

$ws.Range("E28").Value = '  +1.24%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '40.30'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.80%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '10.35'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.40%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.85'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.11%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '155.26'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.39%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.43'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.57%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.17'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.01%  '

$ws.Range("E35").Value = '  +1.80%  '

$ws.Range("E36").Value = '  +3.50%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '18.72'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.46%  '

$ws.Range("E38").Value = '  -0.43%  '

$ws.Range("E39").Value = '  +6.47%  '

$ws.Range("E40").Value = '  -0.46%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '22.51'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -5.29%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.95'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.86%  '

$ws.Range("E43").Value = '  -0.05%  '

$ws.Range("E44").Value = '  +0.07%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.26'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.28%  '

$ws.Range("D46").Value = '2.009.46'
$ws.Range("E46").Value = '  -0.74%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.96'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.70%  '

$ws.Range("D48").Value = '2.831.14'
$ws.Range("E48").Value = '  +1.97%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '83.12'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.32%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '75.14'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.91%  '

$ws.Range("E51").Value = '  +1.95%  '
